# Regenerate the lattice-multiplication exercise table: each of the 15
# cells gets a new "factor x factor" problem (and matching digit-split /
# divider / partial-product placeholder lines), while leaving the table
# structure, fonts, and everything else untouched.
# [char]11 is a vertical-tab, which is how Word represents a manual
# line-break (<w:br/>) inside Range.Text.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$cell.Range.Text = "19 x 65" + [char]11 + "  6    5" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "9|    |"

$cell = $t.Cell(1, 2)
$cell.Range.Text = "86 x 85" + [char]11 + "  8    5" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "6|    |"

$cell = $t.Cell(1, 3)
$cell.Range.Text = "63 x 85" + [char]11 + "  8    5" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "3|    |"

$cell = $t.Cell(2, 1)
$cell.Range.Text = "21 x 54" + [char]11 + "  5    4" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "1|    |"

$cell = $t.Cell(2, 2)
$cell.Range.Text = "54 x 69" + [char]11 + "  6    9" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "4|    |"

$cell = $t.Cell(2, 3)
$cell.Range.Text = "69 x 72" + [char]11 + "  7    2" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "9|    |"

$cell = $t.Cell(3, 1)
$cell.Range.Text = "59 x 17" + [char]11 + "  1    7" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "9|    |"

$cell = $t.Cell(3, 2)
$cell.Range.Text = "45 x 39" + [char]11 + "  3    9" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "5|    |"

$cell = $t.Cell(3, 3)
$cell.Range.Text = "42 x 68" + [char]11 + "  6    8" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "2|    |"

$cell = $t.Cell(4, 1)
$cell.Range.Text = "37 x 47" + [char]11 + "  4    7" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "7|    |"

$cell = $t.Cell(4, 2)
$cell.Range.Text = "80 x 74" + [char]11 + "  7    4" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "0|    |"

$cell = $t.Cell(4, 3)
$cell.Range.Text = "40 x 57" + [char]11 + "  5    7" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "0|    |"

$cell = $t.Cell(5, 1)
$cell.Range.Text = "60 x 36" + [char]11 + "  3    6" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "0|    |"

$cell = $t.Cell(5, 2)
$cell.Range.Text = "31 x 48" + [char]11 + "  4    8" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "1|    |"

$cell = $t.Cell(5, 3)
$cell.Range.Text = "78 x 25" + [char]11 + "  2    5" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "8|    |"
